$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (D = Price, E = Volume) keep their values as plain text
# rather than being auto-converted to numbers/dates by Excel, matching the
# original inlineStr cell type used throughout the sheet.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.830.51'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '2.327.52'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('D5').Value = '302.49'
$ws.Range('D6').Value = '96.35'
$ws.Range('E6').Value = '  +1.51%  '
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.494'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').Value = '34.68'
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('D11').Value = '19.16'
$ws.Range('E11').Value = '  +6.82%  '
$ws.Range('D12').Value = '0.0786'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('D13').Value = '0.120'
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').Value = '6.76'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').Value = '2.683.19'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '2.315.15'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '0.788'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = '42.767.98'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').Value = '12.28'
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('D20').Value = '6.16'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').Value = '0.0₃0891'
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('D22').Value = '68.08'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').Value = '2.30'
$ws.Range('E23').Value = '  +7.00%  '
$ws.Range('D24').Value = '236.16'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('D27').Value = '24.47'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').Value = '166.49'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').Value = '9.14'
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('E31').Value = '  +0.81%  '
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').Value = '5.03'
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('D34').Value = '17.80'
$ws.Range('E34').Value = '  +1.19%  '
$ws.Range('D35').Value = '4.49'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('D36').Value = '0.0702'
$ws.Range('E36').Value = '  +2.96%  '
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('D38').Value = '1.80'
$ws.Range('E38').Value = '  +4.31%  '
$ws.Range('D39').Value = '0.100'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '2.74'
$ws.Range('E40').Value = '  +3.79%  '
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').Value = '20.43'
$ws.Range('E42').Value = '  +15.13%  '
$ws.Range('D43').Value = '1.952.35'
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0280'
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '10.35'
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('E46').Value = '  +3.10%  '
$ws.Range('D47').Value = '2.76'
$ws.Range('E47').Value = '  +1.19%  '
$ws.Range('D48').Value = '2.553.96'
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('D49').Value = '53.50'
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('E50').Value = '  -3.14%  '
$ws.Range('D51').Value = '72.23'
$ws.Range('E51').Value = '  +2.81%  '

# Restore the default (unstyled) look for the price/volume columns so that
# only the text values changed, not the cell formatting/style index.
$ws.Range("D2:E51").Style = "Normal"
